$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose content (columns B:AB) must be swapped.
# Column A (sequential id) stays untouched.
$pairs = @(
    @(8,9),
    @(18,19),
    @(20,21),
    @(22,23),
    @(26,27),
    @(28,29),
    @(32,33),
    @(36,37),
    @(38,39),
    @(43,44),
    @(47,48),
    @(51,52),
    @(56,57),
    @(58,59)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
